$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $c = $ws.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.ClearFormats()
}

    $ws.Range("D2").Value = '27.216.69'
    $ws.Range("E2").Value = '  +0.60%  '
    $ws.Range("D3").Value = '1.896.21'
    $ws.Range("E3").Value = '  +0.15%  '
    $ws.Range("E4").Value = '  +0.00%  '
    Set-TextValue "D5" '307.36'
    $ws.Range("E5").Value = '  +0.10%  '
    Set-TextValue "D6" '1.001'
    $ws.Range("E6").Value = '  +0.05%  '
    Set-TextValue "D7" '0.5198'
    $ws.Range("E7").Value = '  +0.21%  '
    Set-TextValue "D8" '0.3766'
    $ws.Range("E8").Value = '  -0.24%  '
    Set-TextValue "D9" '0.07283'
    $ws.Range("E9").Value = '  +0.87%  '
    Set-TextValue "D10" '21.20'
    $ws.Range("E10").Value = '  +0.40%  '
    Set-TextValue "D11" '0.9005'
    $ws.Range("E11").Value = '  +0.99%  '
    Set-TextValue "D12" '0.08170'
    $ws.Range("E12").Value = '  +6.56%  '
    $ws.Range("B13").Value = 'Litecoin'
    $ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    Set-TextValue "D13" '96.60'
    $ws.Range("E13").Value = '  +2.66%  '
    $ws.Range("B14").Value = 'WrappedEther'
    $ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    $ws.Range("D14").Value = '1.896.70'
    $ws.Range("E14").Value = '  +0.23%  '
    Set-TextValue "D15" '5.279'
    $ws.Range("E15").Value = '  +1.01%  '
    $ws.Range("E16").Value = '  +0.01%  '
    Set-TextValue "D17" '0.000008613'
    $ws.Range("E17").Value = '  +1.23%  '
    Set-TextValue "D18" '14.54'
    $ws.Range("E18").Value = '  +0.36%  '
    Set-TextValue "D19" '1.001'
    $ws.Range("E19").Value = '  +0.03%  '
    $ws.Range("D20").Value = '27.247.98'
    Set-TextValue "D21" '5.085'
    $ws.Range("E21").Value = '  +0.49%  '
    Set-TextValue "D22" '10.69'
    $ws.Range("E22").Value = '  +1.05%  '
    Set-TextValue "D23" '6.398'
    $ws.Range("E23").Value = '  -0.11%  '
    $ws.Range("E24").Value = '  +0.67%  '
    Set-TextValue "D25" '147.29'
    $ws.Range("E25").Value = '  +0.83%  '
    $ws.Range("B26").Value = 'Toncoin'
    $ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    Set-TextValue "D26" '1.744'
    $ws.Range("E26").Value = '  +0.61%  '
    $ws.Range("B27").Value = 'EthereumClassic'
    $ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    Set-TextValue "D27" '18.19'
    $ws.Range("E27").Value = '  +0.78%  '
    Set-TextValue "D28" '115.16'
    $ws.Range("E28").Value = '  +0.60%  '
    Set-TextValue "D29" '4.954'
    $ws.Range("E29").Value = '  -0.25%  '
    Set-TextValue "D30" '4.830'
    $ws.Range("E30").Value = '  +1.00%  '
    Set-TextValue "D31" '0.09236'
    $ws.Range("E31").Value = '  +0.50%  '
    Set-TextValue "D32" '0.7953'
    $ws.Range("E32").Value = '  +2.55%  '
    Set-TextValue "D33" '0.05032'
    $ws.Range("E33").Value = '  -0.29%  '
    Set-TextValue "D34" '1.220'
    $ws.Range("E34").Value = '  -1.20%  '
    Set-TextValue "D35" '3.449'
    Set-TextValue "D36" '2.960'
    $ws.Range("E36").Value = '  -0.49%  '
    Set-TextValue "D37" '2.596'
    $ws.Range("E37").Value = '  +0.39%  '
    Set-TextValue "D38" '0.5673'
    $ws.Range("E38").Value = '  +0.87%  '
    Set-TextValue "D39" '0.01985'
    $ws.Range("E39").Value = '  -0.07%  '
    Set-TextValue "D40" '1.074'
    $ws.Range("E40").Value = '  +0.06%  '
    Set-TextValue "D41" '8.951'
    $ws.Range("E41").Value = '  -0.22%  '
    Set-TextValue "D42" '6.552'
    $ws.Range("E42").Value = '  -1.05%  '
    Set-TextValue "D43" '115.58'
    $ws.Range("E43").Value = '  -2.45%  '
    Set-TextValue "D44" '0.1515'
    $ws.Range("E44").Value = '  -0.05%  '
    Set-TextValue "D45" '0.4873'
    $ws.Range("E45").Value = '  +1.10%  '
    Set-TextValue "D46" '1.001'
    $ws.Range("E46").Value = '  +0.06%  '
    Set-TextValue "D47" '10.10'
    $ws.Range("E47").Value = '  -0.65%  '
    Set-TextValue "D48" '1.623'
    $ws.Range("E48").Value = '  +1.89%  '
    Set-TextValue "D49" '38.20'
    $ws.Range("E49").Value = '  +1.87%  '
    Set-TextValue "D50" '63.46'
    $ws.Range("E50").Value = '  -0.78%  '
    Set-TextValue "D51" '0.05939'
    $ws.Range("E51").Value = '  +0.32%  '
